$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 3275
$ws.Range("I3").Value = 3369
$ws.Range("I4").Value = 794
$ws.Range("I5").Value = 313
$ws.Range("I6").Value = 3816
$ws.Range("I7").Value = 11567

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 105
$ws.Range("I4").Value = 46
$ws.Range("I6").Value = 79
$ws.Range("I7").Value = 384
$ws.Range("I8").Value = 720
$ws.Range("I11").Value = 183
$ws.Range("I13").Value = 19
$ws.Range("I14").Value = 60
$ws.Range("I15").Value = 142
$ws.Range("I20").Value = 294
$ws.Range("I23").Value = 104
$ws.Range("I29").Value = 749
$ws.Range("I30").Value = 40
$ws.Range("I31").Value = 105
$ws.Range("I33").Value = 518
$ws.Range("I34").Value = 52
$ws.Range("I36").Value = 155
$ws.Range("I37").Value = 370
$ws.Range("I38").Value = 8
$ws.Range("I41").Value = 53
$ws.Range("I42").Value = 401
$ws.Range("I43").Value = 103
$ws.Range("I48").Value = 140
$ws.Range("I49").Value = 87
$ws.Range("I51").Value = 110
$ws.Range("I52").Value = 249
$ws.Range("I54").Value = 256
$ws.Range("I60").Value = 60
$ws.Range("I63").Value = 43
$ws.Range("I64").Value = 110
$ws.Range("I65").Value = 256
$ws.Range("I66").Value = 28
$ws.Range("I67").Value = 457
$ws.Range("I70").Value = 22
$ws.Range("I72").Value = 44
$ws.Range("I76").Value = 180
$ws.Range("I78").Value = 157
$ws.Range("I79").Value = 297
$ws.Range("I80").Value = 40
$ws.Range("I86").Value = 67
$ws.Range("I88").Value = 105
$ws.Range("I91").Value = 140
$ws.Range("I95").Value = 184
$ws.Range("I96").Value = 131
$ws.Range("I99").Value = 217
$ws.Range("I101").Value = 11567

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 68
$ws.Range("I3").Value = 86
$ws.Range("I6").Value = 60
$ws.Range("I7").Value = 249

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I2").Value = 81
$ws.Range("I7").Value = 183

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 225
$ws.Range("I6").Value = 234
$ws.Range("I7").Value = 720

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I3").Value = 116
$ws.Range("I7").Value = 384

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I3").Value = 38
$ws.Range("I7").Value = 131

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I2").Value = 17
$ws.Range("I7").Value = 60

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("I2").Value = 11
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I6").Value = 103
$ws.Range("I7").Value = 370

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I2").Value = 60
$ws.Range("I3").Value = 77
$ws.Range("I6").Value = 60
$ws.Range("I7").Value = 217

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I3").Value = 165
$ws.Range("I7").Value = 457

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("I2").Value = 32
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 84
$ws.Range("I5").Value = 12
$ws.Range("I7").Value = 256

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 66
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 184

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I5").Value = 17
$ws.Range("I6").Value = 171
$ws.Range("I7").Value = 518

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 87

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I2").Value = 59
$ws.Range("I3").Value = 49
$ws.Range("I7").Value = 256

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 226
$ws.Range("I4").Value = 33
$ws.Range("I6").Value = 200
$ws.Range("I7").Value = 749

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("I3").Value = 29
$ws.Range("I4").Value = 15
$ws.Range("I6").Value = 77
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item('River North')
$ws.Range("I3").Value = 42
$ws.Range("I7").Value = 180

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("I2").Value = 35
$ws.Range("I7").Value = 79

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I3").Value = 18
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I3").Value = 138
$ws.Range("I4").Value = 35
$ws.Range("I7").Value = 401

$ws = $wb.Worksheets.Item('Boystown')
$ws.Range("I3").Value = 4
$ws.Range("I6").Value = 19

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I6").Value = 64
$ws.Range("I7").Value = 157

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("I3").Value = 36
$ws.Range("I6").Value = 30
$ws.Range("I7").Value = 104

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("I3").Value = 46
$ws.Range("I7").Value = 140

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I3").Value = 94
$ws.Range("I7").Value = 297

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I2").Value = 27
$ws.Range("I3").Value = 35
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I2").Value = 83
$ws.Range("I3").Value = 87
$ws.Range("I7").Value = 294

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I6").Value = 48
$ws.Range("I7").Value = 155

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("I2").Value = 20
$ws.Range("I7").Value = 52

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("I4").Value = 10
$ws.Range("I7").Value = 142

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("I4").Value = 3
$ws.Range("I7").Value = 28

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I3").Value = 37
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 22

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I2").Value = 24
$ws.Range("I3").Value = 37
$ws.Range("I7").Value = 105

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I6").Value = 14
$ws.Range("I7").Value = 67

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I3").Value = 33
$ws.Range("I5").Value = 2
$ws.Range("I7").Value = 110

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("I2").Value = 14
$ws.Range("I7").Value = 60

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I6").Value = 61
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("I2").Value = 8
$ws.Range("I6").Value = 23
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("I2").Value = 8
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("I2").Value = 18
$ws.Range("I7").Value = 46

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range("I5").Value = 3
$ws.Range("I6").Value = 8
